# Update cryptos list - apply latest price/volume snapshot to sheet1
# (GitHub Actions scheduled refresh of coinranking.com data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates for the Price (D) and Volume/1h (E) columns.
# DText marks D values that look numeric (single decimal point) and
# therefore need the cell kept as Text so Excel doesn't silently coerce
# them into a float and drop formatting like trailing zeros.
$updates = @(
    @{ Row = 2;  D = "36.198.11"; E = "  +1.84%  ";  DText = $false },
    @{ Row = 3;  D = "2.023.56";  E = "  +6.87%  ";  DText = $false },
    @{ Row = 5;  D = "243.44";    E = "  -1.10%  ";  DText = $true },
    @{ Row = 6;  D = $null;       E = "  -5.78%  ";  DText = $false },
    @{ Row = 7;  D = $null;       E = "  +0.03%  ";  DText = $false },
    @{ Row = 8;  D = "43.63";     E = "  +1.38%  ";  DText = $true },
    @{ Row = 9;  D = "59.75";     E = "  +4.86%  ";  DText = $true },
    @{ Row = 10; D = "0.354";     E = "  -0.97%  ";  DText = $true },
    @{ Row = 11; D = "0.0708";    E = "  -5.95%  ";  DText = $true },
    @{ Row = 12; D = "0.0978";    E = "  -0.57%  ";  DText = $true },
    @{ Row = 13; D = "14.07";     E = "  -3.42%  ";  DText = $true },
    @{ Row = 14; D = "2.314.54";  E = "  +6.55%  ";  DText = $false },
    @{ Row = 15; D = $null;       E = "  -0.70%  ";  DText = $false },
    @{ Row = 16; D = "2.022.07";  E = "  +6.37%  ";  DText = $false },
    @{ Row = 17; D = $null;       E = "  -4.49%  ";  DText = $false },
    @{ Row = 18; D = "36.218.87"; E = "  +1.91%  ";  DText = $false },
    @{ Row = 19; D = $null;       E = "  -4.42%  ";  DText = $false },
    @{ Row = 20; D = "0.0₃0800";  E = "  -3.61%  ";  DText = $false },
    @{ Row = 21; D = "233.33";    E = "  -5.21%  ";  DText = $true },
    @{ Row = 22; D = "12.36";     E = "  -4.74%  ";  DText = $true },
    @{ Row = 23; D = "4.82";      E = "  -7.14%  ";  DText = $true },
    @{ Row = 24; D = $null;       E = "  +0.02%  ";  DText = $false },
    @{ Row = 25; D = $null;       E = "  -9.36%  ";  DText = $false },
    @{ Row = 26; D = "167.64";    E = "  +0.48%  ";  DText = $true },
    @{ Row = 27; D = "8.62";      E = "  -0.44%  ";  DText = $true },
    @{ Row = 28; D = "19.66";     E = "  +7.04%  ";  DText = $true },
    @{ Row = 29; D = "1.90";      E = "  -11.01%  "; DText = $true },
    @{ Row = 30; D = $null;       E = "  -6.19%  ";  DText = $false },
    @{ Row = 31; D = "21.18";     E = "  +48.75%  "; DText = $true },
    @{ Row = 32; D = $null;       E = "  -2.61%  ";  DText = $false },
    @{ Row = 33; D = "0.0569";    E = "  -6.53%  ";  DText = $true },
    @{ Row = 37; D = $null;       E = "  -7.79%  ";  DText = $false },
    @{ Row = 38; D = "2.12";      E = "  +8.31%  ";  DText = $true },
    @{ Row = 39; D = "0.836";     E = "  -2.15%  ";  DText = $true },
    @{ Row = 40; D = "1.29";      E = "  -12.42%  "; DText = $true },
    @{ Row = 41; D = "94.81";     E = "  -4.34%  ";  DText = $true },
    @{ Row = 42; D = "0.0210";    E = "  -7.88%  ";  DText = $true },
    @{ Row = 43; D = $null;       E = "  +15.79%  "; DText = $false },
    @{ Row = 44; D = $null;       E = "  +0.34%  ";  DText = $false },
    @{ Row = 45; D = "15.36";     E = "  -9.91%  ";  DText = $true },
    @{ Row = 46; D = "1.297.31";  E = "  -1.59%  ";  DText = $false },
    @{ Row = 47; D = $null;       E = "  +0.11%  ";  DText = $false },
    @{ Row = 48; D = "2.77";      E = "  +1.30%  ";  DText = $true },
    @{ Row = 49; D = "2.202.06";  E = $null;         DText = $false },
    @{ Row = 50; D = "2.15";      E = "  -8.38%  ";  DText = $true },
    @{ Row = 51; D = "3.72";      E = "  +12.71%  "; DText = $true }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($r, 4)
        if ($u.DText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}

# Rows 34-36 are reordered: old BinanceUSD/WEMIXToken/Kaspa become
# Kaspa/BinanceUSD/WEMIXToken (with updated price/volume data).
$kaspaCell = $ws.Cells.Item(34, 4)
$kaspaCell.NumberFormat = "@"
$ws.Cells.Item(34, 2).Value = "Kaspa"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$kaspaCell.Value = "0.0882"
$ws.Cells.Item(34, 5).Value = "  +20.16%  "

$busdCell = $ws.Cells.Item(35, 4)
$busdCell.NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = "BinanceUSD"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$busdCell.Value = "1.00"
$ws.Cells.Item(35, 5).Value = "  +0.03%  "

$wemixCell = $ws.Cells.Item(36, 4)
$wemixCell.NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = "WEMIXToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$wemixCell.Value = "1.87"
$ws.Cells.Item(36, 5).Value = "  +0.95%  "
